$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix C column values (Time): row3 0->1, row4 3->2
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2

# D column: convert text "Enabled"/"Disabled" values into boolean FALSE,
# keeping header "Enabled?" in D1 unchanged.
$ws.Range("D2").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("D5").Value = $false

# Update selection to C5
$ws.Range("C5").Select()
